$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the entire content of a paragraph (everything except the
# trailing paragraph mark) with a single run of plain text, then apply a
# character style to that new run. This mirrors what Word does when a user
# selects a paragraph's whole text and types over it with a style applied:
# the many old runs collapse into one clean run.
# ---------------------------------------------------------------------------
function Replace-ParagraphRun($ParaIndex, $NewText, $StyleName) {
    $p = $d.Paragraphs.Item($ParaIndex)
    $full = $p.Range
    # Range covering the paragraph's content but not its end-of-paragraph mark.
    $body = $d.Range($full.Start, $full.End - 1)
    $body.Delete()

    $p2 = $d.Paragraphs.Item($ParaIndex)
    $full2 = $p2.Range
    $insertionPoint = $d.Range($full2.Start, $full2.End - 1)
    $insertionPoint.InsertAfter($NewText)

    $p3 = $d.Paragraphs.Item($ParaIndex)
    $full3 = $p3.Range
    $newRun = $d.Range($full3.Start, $full3.End - 1)
    $newRun.Style = $StyleName
}

# Paragraphs containing the "英仙座：2018年10月30日至..." short-text blurb
# (four occurrences throughout the document) -> single run, style GaNStyle.
$ganStyleText = "英仙座： 2022年1 月 16 日至 25 日、11 月 7 日至 16 日、12 月 6 日至 15 日."
foreach ($idx in @(3, 60, 93, 127)) {
    Replace-ParagraphRun $idx $ganStyleText "GaNStyle"
}

# Paragraph containing the "你现在参加的是全球公益科普活动..." intro text
# -> single run, style GaNParagraph.
$ganParagraphText = "你现在参加的是全球公益科普活动 Globe at Night （全球观星活动），这是一个以观察和记录夜空的可见恒星数来测量你所在地光污染情况的活动。通过定位和观测夜空中的英仙座，并将所肉眼观察到的星等情况与所给出的星等图表作对比，我们可以知道自己社区中的人造光是如何导致光污染的。你所提供数据将和来自全世界的数据一起帮助建立一张全球光污染地图。"
Replace-ParagraphRun 6 $ganParagraphText "GaNParagraph"

# Paragraph containing the "本文檔中的圖表由 Jenik Hollan, CzechGlobe (...)"
# credit line -> single run, style GaNLinks, with the link year updated
# from 2019 to 2022.
$ganLinksText = "本文檔中的圖表由 Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
Replace-ParagraphRun 23 $ganLinksText "GaNLinks"
